$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# 2) Date updated
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# 3) Publisher value was blank -> "Alvearie Team"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# 4) Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# 5) Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> delete it entirely,
#    shifting all subsequent rows up by one.
$ws.Rows.Item(11).Delete()

# 6) "Case Sensitive" row (now row 14 after the deletion) gets value "true" stored as literal text
#    (leading apostrophe forces Excel to keep it as text instead of auto-converting to a boolean).
$ws.Cells.Item(14, 2).Value = "'true"
